$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44/45: Decentraland and EnergySwap swap places (with updated price/volume figures)
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4871"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.65%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "

# Price (D) and Volume(1h) (E) updates for the remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.446.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4769"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3811"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07316"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9308"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07793"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.442"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008815"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.461.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.091"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.009"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.948"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08903"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.321"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.206"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.593"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7512"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.719"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02046"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5543"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05262"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.036"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.615"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1522"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.010"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.666"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06089"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9121"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.59%  "
